$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.524.55'
$ws.Range('E2').Value = '  -7.48%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.547.24'
$ws.Range('E3').Value = '  -1.75%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '295.22'
$ws.Range('E5').Value = '  -4.99%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.59'
$ws.Range('E6').Value = '  -8.17%  '

$ws.Range('E7').Value = '  -4.22%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  -5.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.30'
$ws.Range('E10').Value = '  -8.97%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0799'
$ws.Range('E11').Value = '  -4.27%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.57'
$ws.Range('E12').Value = '  -6.46%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.933.01'
$ws.Range('E13').Value = '  -2.21%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.106'
$ws.Range('E14').Value = '  +0.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.540.29'
$ws.Range('E15').Value = '  -2.40%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.860'
$ws.Range('E16').Value = '  -5.62%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.01'
$ws.Range('E17').Value = '  -5.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.541.97'
$ws.Range('E18').Value = '  -7.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.61'
$ws.Range('E19').Value = '  -1.47%  '

$ws.Range('E20').Value = '  -2.13%  '

$ws.Range('E21').Value = '  -5.26%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.23'
$ws.Range('E22').Value = '  +0.74%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.18'
$ws.Range('E23').Value = '  -6.76%  '

$ws.Range('E24').Value = '  -6.37%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.34'
$ws.Range('E25').Value = '  -1.31%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.09'
$ws.Range('E26').Value = '  -5.03%  '

$ws.Range('E27').Value = '  +0.25%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.91'
$ws.Range('E28').Value = '  -7.25%  '

$ws.Range('E29').Value = '  -4.17%  '

$ws.Range('E30').Value = '  -5.25%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.87'
$ws.Range('E31').Value = '  -5.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.37'
$ws.Range('E32').Value = '  -3.40%  '

$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  -2.30%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.36'
$ws.Range('E34').Value = '  -6.45%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.71'
$ws.Range('E35').Value = '  -3.15%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0785'
$ws.Range('E36').Value = '  -5.86%  '

$ws.Range('E37').Value = '  -8.17%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.16'
$ws.Range('E38').Value = '  +5.23%  '

$ws.Range('E39').Value = '  -3.63%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.53'
$ws.Range('E40').Value = '  -0.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.39'
$ws.Range('E41').Value = '  -5.00%  '

$ws.Range('E42').Value = '  -6.88%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.060.08'
$ws.Range('E43').Value = '  -1.62%  '

$ws.Range('E44').Value = '  -4.50%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.14%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '84.13'
$ws.Range('E46').Value = '  -11.77%  '

$ws.Range('E47').Value = '  +3.36%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.789.18'
$ws.Range('E48').Value = '  -2.25%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.71'
$ws.Range('E49').Value = '  -9.15%  '

$ws.Range('E50').Value = '  -2.98%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.95'
$ws.Range('E51').Value = '  -4.95%  '
